$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ADBS1 ICPL Curncy"
$ws.Range("D31").Value = "NKBS1 ICPL Curncy"
$ws.Range("D30").Value = "NDBS1 ICPL Curncy"
$ws.Range("D27").Value = "KRBS1 TPRA Curncy"

$ws.Range("D14").Select()
